# Apply the "Updated cryptos list" data refresh (Wed May 24 13:53:11 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param([string]$Addr, [string]$NewValue)
    $cell = $ws.Range($Addr)
    # Force text interpretation so numeric-looking strings (e.g. "1.010")
    # keep their exact literal form instead of being parsed as numbers,
    # then restore the default style so no stray number format sticks.
    $cell.NumberFormat = "@"
    $cell.Value = $NewValue
    $cell.Style = "Normal"
}

Set-TextValue 'D2' '26.757.30'
Set-TextValue 'E2' '  -1.91%  '
Set-TextValue 'D3' '1.824.66'
Set-TextValue 'E3' '  -1.32%  '
Set-TextValue 'D4' '1.010'
Set-TextValue 'E4' '  +0.61%  '
Set-TextValue 'D5' '310.83'
Set-TextValue 'E5' '  -0.96%  '
Set-TextValue 'E6' '  +0.55%  '
Set-TextValue 'D7' '0.4575'
Set-TextValue 'E7' '  -0.67%  '
Set-TextValue 'D8' '0.3681'
Set-TextValue 'E8' '  -0.64%  '
Set-TextValue 'D9' '0.07143'
Set-TextValue 'E9' '  -1.94%  '
Set-TextValue 'D10' '0.8739'
Set-TextValue 'E10' '  -1.12%  '
Set-TextValue 'D11' '0.07802'
Set-TextValue 'E11' '  +0.04%  '
Set-TextValue 'D12' '19.38'
Set-TextValue 'E12' '  -2.57%  '
Set-TextValue 'D13' '1.873.29'
Set-TextValue 'E13' '  +0.48%  '
Set-TextValue 'D14' '5.310'
Set-TextValue 'E14' '  -1.16%  '
Set-TextValue 'D15' '6.342'
Set-TextValue 'E15' '  -3.05%  '
Set-TextValue 'D16' '86.87'
Set-TextValue 'E16' '  -5.06%  '
Set-TextValue 'D17' '1.011'
Set-TextValue 'E17' '  +0.62%  '
Set-TextValue 'D18' '0.000008690'
Set-TextValue 'E18' '  -2.81%  '
Set-TextValue 'D19' '1.007'
Set-TextValue 'E19' '  +0.59%  '
Set-TextValue 'D20' '26.827.25'
Set-TextValue 'E20' '  -1.75%  '
Set-TextValue 'D21' '14.40'
Set-TextValue 'E21' '  -2.29%  '
Set-TextValue 'D22' '4.975'
Set-TextValue 'E22' '  -2.76%  '
Set-TextValue 'D23' '2.081.33'
Set-TextValue 'E23' '  +0.98%  '
Set-TextValue 'D24' '10.43'
Set-TextValue 'E24' '  -0.90%  '
Set-TextValue 'D25' '1.993'
Set-TextValue 'E25' '  +3.33%  '
Set-TextValue 'D26' '151.58'
Set-TextValue 'E26' '  -0.03%  '
Set-TextValue 'D27' '18.11'
Set-TextValue 'E27' '  -1.42%  '
Set-TextValue 'D28' '1.991'
Set-TextValue 'E28' '  -2.83%  '
Set-TextValue 'D29' '113.48'
Set-TextValue 'E29' '  -2.04%  '
Set-TextValue 'D30' '4.903'
Set-TextValue 'E30' '  -3.30%  '
Set-TextValue 'D31' '0.08764'
Set-TextValue 'E31' '  -0.74%  '
Set-TextValue 'D32' '3.109'
Set-TextValue 'E32' '  -0.30%  '
Set-TextValue 'D33' '0.7359'
Set-TextValue 'E33' '  -4.26%  '
Set-TextValue 'D34' '4.466'
Set-TextValue 'E34' '  -0.55%  '
Set-TextValue 'D35' '1.122'
Set-TextValue 'E35' '  -3.89%  '
Set-TextValue 'D36' '2.467'
Set-TextValue 'E36' '  -6.67%  '
Set-TextValue 'D37' '1.081'
Set-TextValue 'E37' '  +0.16%  '
Set-TextValue 'D38' '0.01933'
Set-TextValue 'E38' '  -1.14%  '
Set-TextValue 'D39' '0.05108'
Set-TextValue 'E39' '  -2.36%  '
Set-TextValue 'D40' '2.904'
Set-TextValue 'E40' '  -1.82%  '
Set-TextValue 'E41' '  -1.10%  '
Set-TextValue 'D42' '0.4934'
Set-TextValue 'E42' '  -3.71%  '
Set-TextValue 'D43' '0.1587'
Set-TextValue 'E43' '  -2.75%  '
Set-TextValue 'D44' '8.224'
Set-TextValue 'E44' '  -1.78%  '
Set-TextValue 'D45' '1.009'
Set-TextValue 'E45' '  +0.63%  '
Set-TextValue 'D46' '0.4636'
Set-TextValue 'E46' '  -3.38%  '
Set-TextValue 'B47' 'Quant'
Set-TextValue 'C47' 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue 'D47' '102.97'
Set-TextValue 'E47' '  +0.33%  '
Set-TextValue 'B48' 'EnergySwap'
Set-TextValue 'C48' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D48' '10.06'
Set-TextValue 'E48' '  -2.60%  '
Set-TextValue 'D49' '1.584'
Set-TextValue 'E49' '  -3.84%  '
Set-TextValue 'D50' '0.06063'
Set-TextValue 'E50' '  -2.55%  '
Set-TextValue 'D51' '64.52'
Set-TextValue 'E51' '  -1.46%  '

Write-Output "Applied 102 cell updates"
